$d = $word.ActiveDocument

$wdParagraph = 4

# Anchor on the paragraph that must survive: "...LOB1036: ... (Requisito
# fraco)". Expand the found range to the whole paragraph (so its end
# includes the paragraph mark) -- that mark is where the content to be
# removed begins.
$anchor = $d.Content
$anchor.Find.Execute("LOB1036: Geometria Analítica (Requisito fraco)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchor.Expand($wdParagraph)
$deleteStart = $anchor.End

# Anchor on the last bit of text that must be removed (the copyright /
# "Powered by Jekyll" line) and expand to its whole paragraph so the
# deletion also consumes that paragraph's mark.
$tail = $d.Content
$tail.Find.Execute("Creative Commons Attribution", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tail.Expand($wdParagraph)
$deleteEnd = $tail.End

# Delete the empty paragraph + "Ver no Jupiter..." paragraph + the
# copyright paragraph in one shot, leaving the following (already empty)
# paragraph and the page-break paragraph untouched.
$killRange = $d.Range($deleteStart, $deleteEnd)
$killRange.Delete()
